$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so values like "1.00"
# or multi-dot numbers ("63.531.08") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.531.08'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '3.347.68'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '530.89'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').Value = '172.23'
$ws.Range('E6').Value = '  -4.70%  '
$ws.Range('D7').Value = '0.595'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.342.95'
$ws.Range('E8').Value = '  +2.24%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '0.607'
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('D11').Value = '52.95'
$ws.Range('E11').Value = '  -6.98%  '
$ws.Range('E12').Value = '  +2.88%  '
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '9.10'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').Value = '3.894.47'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '3.353.36'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = '17.48'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = '63.533.41'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').Value = '11.18'
$ws.Range('E20').Value = '  +2.58%  '
$ws.Range('D21').Value = '0.963'
$ws.Range('E21').Value = '  +1.89%  '
$ws.Range('D22').Value = '371.59'
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').Value = '11.26'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '81.53'
$ws.Range('E24').Value = '  +2.25%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '3.74'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '4.04'
$ws.Range('E26').Value = '  +6.29%  '
$ws.Range('D27').Value = '6.17'
$ws.Range('E27').Value = '  +2.94%  '
$ws.Range('D28').Value = '2.69'
$ws.Range('E28').Value = '  +3.26%  '
$ws.Range('D29').Value = '11.27'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Value = '8.25'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').Value = '28.79'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = '636.54'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = '6.40'
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('D34').Value = '11.17'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = '57.62'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('D38').Value = '36.40'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').Value = '0.379'
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('D40').Value = '0.0₃0719'
$ws.Range('E40').Value = '  +10.96%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').Value = '2.63'
$ws.Range('E42').Value = '  +6.70%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.946.07'
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '0.124'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +8.25%  '
$ws.Range('D46').Value = '2.69'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('D47').Value = '0.0397'
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.07'
$ws.Range('E48').Value = '  +6.09%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').Value = '2.60'
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').Value = '0.124'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = '136.14'
$ws.Range('E51').Value = '  +4.87%  '
